# Apply updated cryptocurrency price/volume data to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.825.48'
$ws.Range('E2').Value = '  -0.01%  '
$ws.Range('D3').Value = '3.165.73'
$ws.Range('E3').Value = '  -0.94%  '
$ws.Range('D4').Value = '''0.999'
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = '''602.30'
$ws.Range('D6').Value = '''153.92'
$ws.Range('E6').Value = '  +0.58%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').Value = '3.166.50'
$ws.Range('E8').Value = '  -0.76%  '
$ws.Range('E9').Value = '  +2.06%  '
$ws.Range('D10').Value = '''0.157'
$ws.Range('E10').Value = '  -1.16%  '
$ws.Range('D11').Value = '''5.68'
$ws.Range('E11').Value = '  -6.24%  '
$ws.Range('D12').Value = '''0.514'
$ws.Range('E12').Value = '  +0.81%  '
$ws.Range('E13').Value = '  -1.91%  '
$ws.Range('D14').Value = '''38.02'
$ws.Range('E14').Value = '  -3.52%  '
$ws.Range('D15').Value = '3.683.96'
$ws.Range('E15').Value = '  -0.99%  '
$ws.Range('D16').Value = '65.930.21'
$ws.Range('E16').Value = '  -0.02%  '
$ws.Range('E17').Value = '  -1.52%  '
$ws.Range('D18').Value = '3.170.68'
$ws.Range('E19').Value = '  +1.11%  '
$ws.Range('D20').Value = '''506.16'
$ws.Range('E20').Value = '  -0.65%  '
$ws.Range('D21').Value = '''15.27'
$ws.Range('E21').Value = '  -0.77%  '
$ws.Range('E22').Value = '  -1.68%  '
$ws.Range('D23').Value = '''7.95'
$ws.Range('E23').Value = '  -2.08%  '
$ws.Range('E24').Value = '  -3.71%  '
$ws.Range('D25').Value = '''84.22'
$ws.Range('E25').Value = '  -0.62%  '
$ws.Range('D26').Value = '''0.999'
$ws.Range('E26').Value = '  +0.03%  '
$ws.Range('D27').Value = '''2.98'
$ws.Range('E27').Value = '  -0.72%  '
$ws.Range('E28').Value = '  -1.68%  '
$ws.Range('D29').Value = '''2.37'
$ws.Range('E29').Value = '  +4.36%  '
$ws.Range('D30').Value = '''7.16'
$ws.Range('E30').Value = '  +4.53%  '
$ws.Range('E31').Value = '  +4.60%  '
$ws.Range('D32').Value = '''27.81'
$ws.Range('E32').Value = '  -0.79%  '
$ws.Range('E33').Value = '  +0.20%  '
$ws.Range('D34').Value = '''1.18'
$ws.Range('E34').Value = '  -3.33%  '
$ws.Range('E35').Value = '  -1.31%  '
$ws.Range('D36').Value = '''506.96'
$ws.Range('E36').Value = '  +4.62%  '
$ws.Range('D37').Value = '''55.23'
$ws.Range('E37').Value = '  +0.45%  '
$ws.Range('E38').Value = '  -3.07%  '
$ws.Range('E39').Value = '  -0.53%  '
$ws.Range('E40').Value = '  +5.47%  '
$ws.Range('D41').Value = '''8.73'
$ws.Range('E41').Value = '  -1.58%  '
$ws.Range('D42').Value = '0.0₃0674'
$ws.Range('E42').Value = '  +5.14%  '
$ws.Range('D43').Value = '''2.84'
$ws.Range('E43').Value = '  -2.61%  '
$ws.Range('E44').Value = '  -2.48%  '
$ws.Range('E45').Value = '  -0.09%  '
$ws.Range('D46').Value = '2.818.36'
$ws.Range('E46').Value = '  -4.28%  '
$ws.Range('D47').Value = '''27.74'
$ws.Range('E47').Value = '  -2.31%  '
$ws.Range('E48').Value = '  -0.07%  '
$ws.Range('E49').Value = '  +2.16%  '
$ws.Range('E50').Value = '  +0.15%  '
$ws.Range('D51').Value = '''2.59'
$ws.Range('E51').Value = '  +0.32%  '
